# NATMI LR-pairs re-run with updated TPM values: existing FAPs->* rows get
# refreshed specificity/weight numbers, and a brand-new "MuSCs" sending
# cluster block (rows 7-11, mirroring the same 5 target clusters) is added.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns: A Sending cluster, B Ligand symbol, C Receptor symbol, D Target cluster,
# E..T the various expression / specificity metrics.
$headerRow = @(
    "Sending cluster","Ligand symbol","Receptor symbol","Target cluster",
    "Ligand-expressing cells","Ligand detection rate","Ligand average expression value",
    "Ligand total expression value","Ligand derived specificity of average expression value",
    "Ligand derived specificity of total expression value","Receptor-expressing cells",
    "Receptor detection rate","Receptor average expression value","Receptor total expression value",
    "Receptor derived specificity of average expression value","Receptor derived specificity of total expression value",
    "Edge average expression weight","Edge total expression weight",
    "Edge average expression derived specificity","Edge total expression derived specificity"
)
for ($c = 0; $c -lt $headerRow.Count; $c++) {
    $ws.Cells.Item(1, $c + 1).Value = $headerRow[$c]
}

$rows = @(
    @("FAPs","Lama1","Itga2","ECs",3,1,0.5587383333333333,1.676215,0.8486764927018626,0.8937587278261895,3,1,3.031236,9.093708000000001,0.6569357730026921,0.6780946256479073,1.69366775058,15.24300975522,0.5575259477623117,0.6060529899648499),
    @("FAPs","Lama1","Itga2","FAPs",3,1,0.5587383333333333,1.676215,0.8486764927018626,0.8937587278261895,3,1,1.067915666666667,3.203747,0.2314409052885859,0.2388952474211406,0.5966854197338889,5.370168777605,0.196418455768061,0.2135147124188414),
    @("FAPs","Lama1","Itga2","Inflammatory-Mac",3,1,0.5587383333333333,1.676215,0.8486764927018626,0.8937587278261895,1,0.3333333333333333,0.04253966666666667,0.127619,0.009219285072143351,0.009516223528461688,0.02376854245388889,0.213916882085,0.007824190520245258,0.008505207834507571),
    @("FAPs","Lama1","Itga2","MuSCs",3,1,0.5587383333333333,1.676215,0.8486764927018626,0.8937587278261895,2,1,0.4319364999999999,0.8638729999999999,0.09361017700884301,0.06441680759293508,0.2413394801158333,1.448036880695,0.07944475670506543,0.05757308400488608),
    @("FAPs","Lama1","Itga2","Resolving-Mac",3,1,0.5587383333333333,1.676215,0.8486764927018626,0.8937587278261895,2,0.6666666666666666,0.04057666666666667,0.12173,0.008793859627735762,0.009077095809555328,0.02267173910555556,0.20404565195,0.007463141946179294,0.008112733603104605),
    @("MuSCs","Lama1","Itga2","ECs",2,1,0.09962599999999999,0.199252,0.1513235072981373,0.1062412721738106,3,1,3.031236,9.093708000000001,0.6569357730026921,0.6780946256479073,0.301989917736,1.811939506416,0.09940982524038038,0.07204163568305752),
    @("MuSCs","Lama1","Itga2","FAPs",2,1,0.09962599999999999,0.199252,0.1513235072981373,0.1062412721738106,3,1,1.067915666666667,3.203747,0.2314409052885859,0.2388952474211406,0.1063921662073333,0.638352997244,0.03502244952052484,0.02538053500229922),
    @("MuSCs","Lama1","Itga2","Inflammatory-Mac",2,1,0.09962599999999999,0.199252,0.1513235072981373,0.1062412721738106,1,0.3333333333333333,0.04253966666666667,0.127619,0.009219285072143351,0.009516223528461688,0.004238056831333333,0.025428340988,0.001395094551898093,0.001011015693954118),
    @("MuSCs","Lama1","Itga2","MuSCs",2,1,0.09962599999999999,0.199252,0.1513235072981373,0.1062412721738106,2,1,0.4319364999999999,0.8638729999999999,0.09361017700884301,0.06441680759293508,0.04303210574899999,0.172128422996,0.01416542030377758,0.006843723588049004),
    @("MuSCs","Lama1","Itga2","Resolving-Mac",2,1,0.09962599999999999,0.199252,0.1513235072981373,0.1062412721738106,2,0.6666666666666666,0.04057666666666667,0.12173,0.008793859627735762,0.009077095809555328,0.004042490993333333,0.02425494596,0.001330717681556468,0.000964362206450723)
)

for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $i + 2
    $values = $rows[$i]
    for ($c = 0; $c -lt $values.Count; $c++) {
        $ws.Cells.Item($r, $c + 1).Value = $values[$c]
    }
}
